$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = " UserId"
$ws.Range("B1").Value = " UserName"
$ws.Range("C1").Value = "FirstName"
$ws.Range("D1").Value = "LastName"
$ws.Range("E1").Value = "Email"
$ws.Range("F1").Value = "Password"
$ws.Range("G1").Value = "Phone"
